$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 89 ---
$ws.Cells.Item(89, 1).Value2 = "WAT96"
$ws.Cells.Item(89, 2).Value2 = "WAT-305"
$ws.Cells.Item(89, 3).Value2 = "Verify that system provides the filter option ""Filter by organization"" in the Author search result page"
$ws.Cells.Item(89, 4).Value2 = "Y"

# --- Row 90 ---
$ws.Cells.Item(90, 1).Value2 = "WAT97"
$ws.Cells.Item(90, 2).Value2 = "WAT-304"
$ws.Cells.Item(90, 3).Value2 = "Verify that system provides the filter option ""Filter by author name"" in the Author search result page"
$ws.Cells.Item(90, 4).Value2 = "Y"

# Carry over the same formatting used by the previous data row (row 88)
# onto the two newly-added rows.
[void]$ws.Range("A88:E88").Copy()
[void]$ws.Range("A89:E90").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the saved selection state recorded for this sheet after the edit.
[void]$ws.Range("C89:C90").Select()
